# Implement basket-based elective scheduling with common time slots across all branches
# Updates the Section_A and Section_B timetable sheets with the new slot assignments.

$wb = $excel.ActiveWorkbook

# --- Section_A (sheet 1) ---
$wsA = $wb.Worksheets.Item("Section_A")

$wsA.Range("C2").Value = "Free"
$wsA.Range("D2").Value = "CS304"
$wsA.Range("E2").Value = "CS307"
$wsA.Range("F2").Value = "MA261"

$wsA.Range("C3").Value = "MA261"
$wsA.Range("D3").Value = "MA262"
$wsA.Range("E3").Value = "CS304"
$wsA.Range("F3").Value = "CS304"

$wsA.Range("B5").Value = "MA262"
$wsA.Range("C5").Value = "CS307"
$wsA.Range("D5").Value = "DA261"

$wsA.Range("F6").Value = "CS307 (Tutorial)"

$wsA.Range("B7").Value = "DA262"
$wsA.Range("C7").Value = "DA262"
$wsA.Range("F7").Value = "DA261"

# --- Section_B (sheet 2) ---
$wsB = $wb.Worksheets.Item("Section_B")

$wsB.Range("C2").Value = "Free"
$wsB.Range("D2").Value = "MA262"
$wsB.Range("E2").Value = "DA262"
$wsB.Range("F2").Value = "MA261"

$wsB.Range("C3").Value = "DA262"
$wsB.Range("D3").Value = "MA261"
$wsB.Range("E3").Value = "CS304"
$wsB.Range("F3").Value = "DA261"

$wsB.Range("B5").Value = "CS304"
$wsB.Range("C5").Value = "CS307"
$wsB.Range("D5").Value = "CS304"

$wsB.Range("B7").Value = "MA262"
$wsB.Range("D7").Value = "DA261"
$wsB.Range("F7").Value = "DA262"

$wsB.Range("D8").Value = "CS307 (Tutorial)"
